# "More corrections to start times"
#
# Corrects several event start times (column F, time-of-day values) and
# one audience figure (H10) on Sheet1, and updates the active selection
# to match the author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Corrected start times (column F) ---
$ws.Range("F3").Value  = 0.54166666666666663   # 16:15 -> 13:00
$ws.Range("F7").Value  = 0.54166666666666663   # 15:30 -> 13:00
$ws.Range("F8").Value  = 0.5625                # 14:30 -> 13:30
$ws.Range("F10").Value = 0.63541666666666663   # 17:15 -> 15:15
$ws.Range("F13").Value = 0.625                 # 16:00 -> 15:00

# --- Corrected audience figure for the Grand National (row 10) ---
$ws.Range("H10").Value = 4.6399999999999997    # was 11.2

# --- Scroll/select to match the saved view state ---
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 3
$ws.Range("F21").Select()
